# Scheduled runner: refresh cached market-price / profit figures across all sheets.
# Applies updated currentAveragePrice / LevePrice / LeveProfit values (columns H-N)
# for the affected leve rows on each class sheet.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H39").Value2 = 1777.625
$ws.Range("I39").Value2 = 2086.8
$ws.Range("J39").Value2 = 1262.3334
$ws.Range("K39").Value2 = 6260.400000000001
$ws.Range("L39").Value2 = 3787.0002
$ws.Range("M39").Value2 = -5964.400000000001
$ws.Range("N39").Value2 = -4379.0002

$ws.Range("H40").Value2 = 4788.35
$ws.Range("I40").Value2 = 4521.625
$ws.Range("J40").Value2 = 4966.1665
$ws.Range("K40").Value2 = 4521.625
$ws.Range("L40").Value2 = 4966.1665
$ws.Range("M40").Value2 = -4346.625
$ws.Range("N40").Value2 = -5316.1665

$ws.Range("H138").Value2 = 2308.7627
$ws.Range("I138").Value2 = 1632.3684
$ws.Range("J138").Value2 = 2630.05
$ws.Range("K138").Value2 = 4897.1052
$ws.Range("L138").Value2 = 7890.150000000001
$ws.Range("M138").Value2 = 242.8948
$ws.Range("N138").Value2 = -18170.15

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H102").Value2 = 6053.375
$ws.Range("I102").Value2 = 6089.5713
$ws.Range("K102").Value2 = 6089.5713
$ws.Range("M102").Value2 = -4467.5713

$ws.Range("H128").Value2 = 99990
$ws.Range("J128").Value2 = 99990
$ws.Range("L128").Value2 = 99990
$ws.Range("N128").Value2 = -109950

$ws.Range("H132").Value2 = 1745.2069
$ws.Range("I132").Value2 = 1496.42
$ws.Range("K132").Value2 = 4489.26
$ws.Range("M132").Value2 = -1959.26

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H8").Value2 = 220.8
$ws.Range("I8").Value2 = 754
$ws.Range("J8").Value2 = 87.5
$ws.Range("K8").Value2 = 754
$ws.Range("L8").Value2 = 87.5
$ws.Range("M8").Value2 = -614
$ws.Range("N8").Value2 = -367.5

$ws.Range("H86").Value2 = 2335.3333
$ws.Range("I86").Value2 = 2236
$ws.Range("K86").Value2 = 2236
$ws.Range("M86").Value2 = -1113

$ws.Range("H89").Value2 = 2335.3333
$ws.Range("I89").Value2 = 2236
$ws.Range("K89").Value2 = 11180
$ws.Range("M89").Value2 = -5564

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value2 = 2941.95
$ws.Range("J31").Value2 = 6424.8335
$ws.Range("L31").Value2 = 6424.8335
$ws.Range("N31").Value2 = -7014.8335

$ws.Range("H34").Value2 = 2941.95
$ws.Range("J34").Value2 = 6424.8335
$ws.Range("L34").Value2 = 6424.8335
$ws.Range("N34").Value2 = -6828.8335

$ws.Range("H58").Value2 = 1977.9259
$ws.Range("I58").Value2 = 1605.0667
$ws.Range("K58").Value2 = 1605.0667
$ws.Range("M58").Value2 = -1402.0667

$ws.Range("H68").Value2 = 20663.334
$ws.Range("J68").Value2 = 20663.334
$ws.Range("L68").Value2 = 20663.334
$ws.Range("N68").Value2 = -22161.334

$ws.Range("H71").Value2 = 20663.334
$ws.Range("J71").Value2 = 20663.334
$ws.Range("L71").Value2 = 61990.00199999999
$ws.Range("N71").Value2 = -69478.00199999999

$ws.Range("H86").Value2 = 2238.96
$ws.Range("I86").Value2 = 1840.421
$ws.Range("J86").Value2 = 3501
$ws.Range("K86").Value2 = 1840.421
$ws.Range("L86").Value2 = 3501
$ws.Range("M86").Value2 = -717.421
$ws.Range("N86").Value2 = -5747

$ws.Range("H89").Value2 = 2238.96
$ws.Range("I89").Value2 = 1840.421
$ws.Range("J89").Value2 = 3501
$ws.Range("K89").Value2 = 9202.105
$ws.Range("L89").Value2 = 17505
$ws.Range("M89").Value2 = -3586.105
$ws.Range("N89").Value2 = -28737

$ws.Range("H97").Value2 = 21913
$ws.Range("J97").Value2 = 21913
$ws.Range("L97").Value2 = 21913
$ws.Range("N97").Value2 = -23895

$ws.Range("H99").Value2 = 3605.524
$ws.Range("I99").Value2 = 2367.8
$ws.Range("J99").Value2 = 4730.727
$ws.Range("K99").Value2 = 2367.8
$ws.Range("L99").Value2 = 4730.727
$ws.Range("M99").Value2 = -869.8000000000002
$ws.Range("N99").Value2 = -7726.727

$ws.Range("H105").Value2 = 1579.3462
$ws.Range("I105").Value2 = 1390.1904
$ws.Range("K105").Value2 = 1390.1904
$ws.Range("M105").Value2 = 356.8096

$ws.Range("H126").Value2 = 3605.524
$ws.Range("I126").Value2 = 2367.8
$ws.Range("J126").Value2 = 4730.727
$ws.Range("K126").Value2 = 7103.400000000001
$ws.Range("L126").Value2 = 14192.181
$ws.Range("M126").Value2 = -4633.400000000001
$ws.Range("N126").Value2 = -19132.181

$ws.Range("H132").Value2 = 3910.7778
$ws.Range("I132").Value2 = 3910.7778
$ws.Range("K132").Value2 = 11732.3334
$ws.Range("M132").Value2 = -9202.3334

$ws.Range("H134").Value2 = 2419.6274
$ws.Range("I134").Value2 = 1871.4667
$ws.Range("J134").Value2 = 6530.8335
$ws.Range("K134").Value2 = 5614.4001
$ws.Range("L134").Value2 = 19592.5005
$ws.Range("M134").Value2 = -3079.4001
$ws.Range("N134").Value2 = -24662.5005

$ws.Range("H136").Value2 = 1977.9259
$ws.Range("I136").Value2 = 1605.0667
$ws.Range("K136").Value2 = 4815.2001
$ws.Range("M136").Value2 = -2265.2001

$ws.Range("H141").Value2 = 134261.8
$ws.Range("J141").Value2 = 134261.8
$ws.Range("L141").Value2 = 134261.8
$ws.Range("N141").Value2 = -144621.8

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H55").Value2 = 9695629
$ws.Range("I55").Value2 = 145514.14
$ws.Range("J55").Value2 = 20837430
$ws.Range("K55").Value2 = 436542.42
$ws.Range("L55").Value2 = 62512290
$ws.Range("M55").Value2 = -436365.42
$ws.Range("N55").Value2 = -62512644

$ws.Range("H99").Value2 = 13846
$ws.Range("I99").Value2 = 1612.5
$ws.Range("J99").Value2 = 18739.4
$ws.Range("K99").Value2 = 4837.5
$ws.Range("L99").Value2 = 56218.2
$ws.Range("M99").Value2 = -2591.5
$ws.Range("N99").Value2 = -60710.2

$ws.Range("H105").Value2 = 13532.667
$ws.Range("J105").Value2 = 13532.667
$ws.Range("L105").Value2 = 40598.001
$ws.Range("N105").Value2 = -45840.001

$ws.Range("H113").Value2 = 1924.5652
$ws.Range("J113").Value2 = 1993.8636
$ws.Range("L113").Value2 = 5981.5908
$ws.Range("N113").Value2 = -10321.5908

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H20").Value2 = 22081.666
$ws.Range("I20").Value2 = 24000
$ws.Range("J20").Value2 = 18245
$ws.Range("K20").Value2 = 24000
$ws.Range("L20").Value2 = 18245
$ws.Range("M20").Value2 = -23755
$ws.Range("N20").Value2 = -18735

$ws.Range("H24").Value2 = 7695784.5
$ws.Range("I24").Value2 = 9093609
$ws.Range("J24").Value2 = 7750
$ws.Range("K24").Value2 = 9093609
$ws.Range("L24").Value2 = 7750
$ws.Range("M24").Value2 = -9093436
$ws.Range("N24").Value2 = -8096

$ws.Range("H70").Value2 = 5881.5454
$ws.Range("I70").Value2 = 5944.1333
$ws.Range("J70").Value2 = 5747.4287
$ws.Range("K70").Value2 = 5944.1333
$ws.Range("L70").Value2 = 5747.4287
$ws.Range("M70").Value2 = -5674.1333
$ws.Range("N70").Value2 = -6287.4287

$ws.Range("H73").Value2 = 5881.5454
$ws.Range("I73").Value2 = 5944.1333
$ws.Range("J73").Value2 = 5747.4287
$ws.Range("K73").Value2 = 5944.1333
$ws.Range("L73").Value2 = 5747.4287
$ws.Range("M73").Value2 = -5008.1333
$ws.Range("N73").Value2 = -7619.4287

$ws.Range("H122").Value2 = 4110.4
$ws.Range("I122").Value2 = 4250.6665
$ws.Range("K122").Value2 = 12751.9995
$ws.Range("M122").Value2 = -10301.9995

$ws.Range("H126").Value2 = 1837.2727
$ws.Range("I126").Value2 = 1826.625
$ws.Range("J126").Value2 = 1865.6666
$ws.Range("K126").Value2 = 5479.875
$ws.Range("L126").Value2 = 5596.9998
$ws.Range("M126").Value2 = -3009.875
$ws.Range("N126").Value2 = -10536.9998

$ws.Range("H132").Value2 = 3294.0483
$ws.Range("I132").Value2 = 2246.0715
$ws.Range("K132").Value2 = 6738.2145
$ws.Range("M132").Value2 = -4208.2145

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value2 = 29992.58
$ws.Range("I7").Value2 = 24171.916
$ws.Range("K7").Value2 = 24171.916
$ws.Range("M7").Value2 = -24059.916

$ws.Range("H22").Value2 = 1724.5
$ws.Range("I22").Value2 = 1466
$ws.Range("K22").Value2 = 1466
$ws.Range("M22").Value2 = -1171

$ws.Range("H27").Value2 = 1724.5
$ws.Range("I27").Value2 = 1466
$ws.Range("K27").Value2 = 1466
$ws.Range("M27").Value2 = -1359

$ws.Range("H40").Value2 = 16645.139
$ws.Range("I40").Value2 = 8535.700000000001
$ws.Range("K40").Value2 = 8535.700000000001
$ws.Range("M40").Value2 = -8399.700000000001

$ws.Range("H96").Value2 = 0
$ws.Range("J96").Value2 = 0
$ws.Range("L96").Value2 = 0
$ws.Range("N96").ClearContents()

$ws.Range("H122").Value2 = 95212.55
$ws.Range("I122").Value2 = 146465.72
$ws.Range("K122").Value2 = 439397.16
$ws.Range("M122").Value2 = -436947.16

$ws.Range("H126").Value2 = 29992.58
$ws.Range("I126").Value2 = 24171.916
$ws.Range("K126").Value2 = 72515.74800000001
$ws.Range("M126").Value2 = -70045.74800000001

$ws.Range("H132").Value2 = 3846.7021
$ws.Range("J132").Value2 = 5289.75
$ws.Range("L132").Value2 = 15869.25
$ws.Range("N132").Value2 = -20929.25

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H51").Value2 = 1100
$ws.Range("I51").Value2 = 1100
$ws.Range("J51").Value2 = 0
$ws.Range("K51").Value2 = 1100
$ws.Range("L51").Value2 = 0
$ws.Range("M51").Value2 = -590
$ws.Range("N51").ClearContents()

$ws.Range("H52").Value2 = 0
$ws.Range("J52").Value2 = 0
$ws.Range("L52").Value2 = 0
$ws.Range("N52").ClearContents()

$ws.Range("H123").Value2 = 0
$ws.Range("J123").Value2 = 0
$ws.Range("L123").Value2 = 0
$ws.Range("N123").ClearContents()

$ws.Range("H125").Value2 = 60882.5
$ws.Range("J125").Value2 = 60882.5
$ws.Range("L125").Value2 = 60882.5
$ws.Range("N125").Value2 = -70722.5
